$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Falling Creek Reservoir (row 2): add DOC value
$ws.Range("H2").Value = 2.59

# Lake Mendota (row 3): add DOC and DIC values, update data provider citation
$ws.Range("H3").Value = 4.87
$ws.Range("I3").Value = 45.16
$ws.Range("N3").Value = "LTER- North Temperate Lakes; Hart et al. 2017"

# Toolik Lake (row 5): update data provider citation
$ws.Range("N5").Value = "NEON Relocatable Aquatic; LTER- Arctic; Kling et al. 2000"

# Update header for "Data providers" column (N1)
$ws.Range("N1").Value = "Data providers"
